# Equity curve update — v0.0.15 + v2.2.9 Up_and_Down
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk update columns A (Equity) and B (DrawdownPct) for data rows 2-117 ---
$data = New-Object 'object[,]' 116,2

$data[0,0] = 1000
$data[0,1] = 0
$data[1,0] = 1000
$data[1,1] = 0
$data[2,0] = 1000
$data[2,1] = 0
$data[3,0] = 1000
$data[3,1] = 0
$data[4,0] = 1000
$data[4,1] = 0
$data[5,0] = 1000
$data[5,1] = 0
$data[6,0] = 1000
$data[6,1] = 0
$data[7,0] = 1000
$data[7,1] = 0
$data[8,0] = 995.0329786537201
$data[8,1] = 0.00496702134627991
$data[9,0] = 1105.74449373372
$data[9,1] = 0
$data[10,0] = 1174.98290685372
$data[10,1] = 0
$data[11,0] = 1163.36545365372
$data[11,1] = 0.009887338047417704
$data[12,0] = 1259.37492761372
$data[12,1] = 0
$data[13,0] = 1369.28380665372
$data[13,1] = 0
$data[14,0] = 1337.86442457372
$data[14,1] = 0.02294585090930357
$data[15,0] = 1340.42776805372
$data[15,1] = 0.02107381863407765
$data[16,0] = 1194.525367683652
$data[16,1] = 0.1276276241060247
$data[17,0] = 1194.525367683652
$data[17,1] = 0.1276276241060247
$data[18,0] = 1194.525367683652
$data[18,1] = 0.1276276241060247
$data[19,0] = 1194.525367683652
$data[19,1] = 0.1276276241060247
$data[20,0] = 1194.525367683652
$data[20,1] = 0.1276276241060247
$data[21,0] = 1343.729449920652
$data[21,1] = 0.01866257134488269
$data[22,0] = 1365.781839180652
$data[22,1] = 0.002557517627865491
$data[23,0] = 1402.300164880652
$data[23,1] = 0
$data[24,0] = 1498.643234240652
$data[24,1] = 0
$data[25,0] = 1700.812859800652
$data[25,1] = 0
$data[26,0] = 1671.963466660653
$data[26,1] = 0.01696212077287629
$data[27,0] = 1833.096896680653
$data[27,1] = 0
$data[28,0] = 1671.554434480652
$data[28,1] = 0.08812543542707385
$data[29,0] = 1794.083979660652
$data[29,1] = 0.0212825176293977
$data[30,0] = 1844.366322380652
$data[30,1] = 0
$data[31,0] = 1776.212808240653
$data[31,1] = 0.03695226556296549
$data[32,0] = 1768.232472580652
$data[32,1] = 0.04127913683748508
$data[33,0] = 1738.213213740652
$data[33,1] = 0.05755532800175012
$data[34,0] = 1730.311991300652
$data[34,1] = 0.06183930474981902
$data[35,0] = 1769.589180140652
$data[35,1] = 0.04054354134133187
$data[36,0] = 1753.806934380652
$data[36,1] = 0.04910054304348221
$data[37,0] = 1875.430885680652
$data[37,1] = 0
$data[38,0] = 2067.290543740653
$data[38,1] = 0
$data[39,0] = 2170.024951320652
$data[39,1] = 0
$data[40,0] = 2208.625469640652
$data[40,1] = 0
$data[41,0] = 2164.377614020652
$data[41,1] = 0.02003411453332515
$data[42,0] = 2174.987201800653
$data[42,1] = 0.01523040837044798
$data[43,0] = 2105.722088421733
$data[43,1] = 0.04659159401782242
$data[44,0] = 2105.722088421733
$data[44,1] = 0.04659159401782242
$data[45,0] = 2105.722088421733
$data[45,1] = 0.04659159401782242
$data[46,0] = 2091.493446887773
$data[46,1] = 0.05303390020759702
$data[47,0] = 2070.857365117673
$data[47,1] = 0.06237730498752003
$data[48,0] = 2199.073040497673
$data[48,1] = 0.004325056137532202
$data[49,0] = 2223.555573537673
$data[49,1] = 0
$data[50,0] = 2502.058077517673
$data[50,1] = 0
$data[51,0] = 2635.801124657673
$data[51,1] = 0
$data[52,0] = 2948.605095857673
$data[52,1] = 0
$data[53,0] = 3036.572085657673
$data[53,1] = 0
$data[54,0] = 3485.156272717672
$data[54,1] = 0
$data[55,0] = 3443.273836237673
$data[55,1] = 0.01201737689866633
$data[56,0] = 3657.652986457673
$data[56,1] = 0
$data[57,0] = 3623.990416237673
$data[57,1] = 0.009203325286634478
$data[58,0] = 4405.366923497673
$data[58,1] = 0
$data[59,0] = 4921.256171377672
$data[59,1] = 0
$data[60,0] = 6147.644518757674
$data[60,1] = 0
$data[61,0] = 7087.730342497673
$data[61,1] = 0
$data[62,0] = 6663.975518277673
$data[62,1] = 0.05978709738421994
$data[63,0] = 6012.536057457673
$data[63,1] = 0.1516979672030112
$data[64,0] = 6164.608146817673
$data[64,1] = 0.1302422850577435
$data[65,0] = 7205.592235637674
$data[65,1] = 0
$data[66,0] = 8991.236333837674
$data[66,1] = 0
$data[67,0] = 10603.22457659767
$data[67,1] = 0
$data[68,0] = 8362.903039377674
$data[68,1] = 0.2112868138400655
$data[69,0] = 9428.234580157674
$data[69,1] = 0.1108144025387632
$data[70,0] = 10887.94263567767
$data[70,1] = 0
$data[71,0] = 10592.81785717767
$data[71,1] = 0.02710565148763133
$data[72,0] = 10305.50952711767
$data[72,1] = 0.05349340348758636
$data[73,0] = 10748.06070107767
$data[73,1] = 0.0128474165671697
$data[74,0] = 11076.71296871767
$data[74,1] = 0
$data[75,0] = 10494.94941648466
$data[75,1] = 0.05252131691739259
$data[76,0] = 10378.08303972466
$data[76,1] = 0.06307195383378172
$data[77,0] = 10502.01428828466
$data[77,1] = 0.05188350389290131
$data[78,0] = 10557.54710882429
$data[78,1] = 0.04687002916475191
$data[79,0] = 9954.351817215951
$data[79,1] = 0.101326192587227
$data[80,0] = 9353.833724715951
$data[80,1] = 0.1555406598390151
$data[81,0] = 9404.11482859595
$data[81,1] = 0.1510013074136158
$data[82,0] = 9412.02093835595
$data[82,1] = 0.1502875478549518
$data[83,0] = 9664.619004606591
$data[83,1] = 0.1274831232062301
$data[84,0] = 9307.834799078832
$data[84,1] = 0.1596934193956657
$data[85,0] = 9214.23462285883
$data[85,1] = 0.1681435955882188
$data[86,0] = 9275.208612428831
$data[86,1] = 0.1626388948938702
$data[87,0] = 9168.324954508831
$data[87,1] = 0.1722882970424909
$data[88,0] = 8910.308145478832
$data[88,1] = 0.1955819230269032
$data[89,0] = 9285.039637138831
$data[89,1] = 0.1617513549948257
$data[90,0] = 9873.362009658211
$data[90,1] = 0.1086379111256118
$data[91,0] = 10491.63365630979
$data[91,1] = 0.05282066205563274
$data[92,0] = 11061.20575170199
$data[92,1] = 0.001399983646725889
$data[93,0] = 11514.49190964549
$data[93,1] = 0
$data[94,0] = 11409.49776525739
$data[94,1] = 0.009118434856873359
$data[95,0] = 12068.98083390739
$data[95,1] = 0
$data[96,0] = 10634.3414653148
$data[96,1] = 0.1188699682546532
$data[97,0] = 10797.2394533748
$data[97,1] = 0.105372723516113
$data[98,0] = 10250.7891532748
$data[98,1] = 0.1506499766346834
$data[99,0] = 10925.5562169348
$data[99,1] = 0.09474077659980851
$data[100,0] = 11922.12494747941
$data[100,1] = 0.01216804371876978
$data[101,0] = 13106.48575918029
$data[101,1] = 0
$data[102,0] = 12977.94980624149
$data[102,1] = 0.00980704937238952
$data[103,0] = 13068.67571802653
$data[103,1] = 0.002884834413166448
$data[104,0] = 13473.13489587653
$data[104,1] = 0
$data[105,0] = 13933.27357752653
$data[105,1] = 0
$data[106,0] = 12841.66498559055
$data[106,1] = 0.0783454502534624
$data[107,0] = 12732.59379549055
$data[107,1] = 0.08617355967032791
$data[108,0] = 12094.70699474055
$data[108,1] = 0.1319551053495049
$data[109,0] = 12147.94714979055
$data[109,1] = 0.1281340252024911
$data[110,0] = 11874.87892324055
$data[110,1] = 0.1477323073312823
$data[111,0] = 12206.37068414055
$data[111,1] = 0.1239409305915997
$data[112,0] = 11854.83784835698
$data[112,1] = 0.1491706681566871
$data[113,0] = 11195.68132401251
$data[113,1] = 0.1964787555689415
$data[114,0] = 11342.21356593251
$data[114,1] = 0.1859620423855911
$data[115,0] = 11342.21356593251
$data[115,1] = 0.1859620423855911

$ws.Range("A2:B117").Value = $data

# --- DrawdownDuration (column C) cell additions: new drawdown-duration markers ---
$ws.Range("C11").Value = 14
$ws.Range("C11").NumberFormat = "0"

$ws.Range("C14").Value = 14
$ws.Range("C14").NumberFormat = "0"

$ws.Range("C25").Value = 70
$ws.Range("C25").NumberFormat = "0"

$ws.Range("C51").Value = 63
$ws.Range("C51").NumberFormat = "0"

$ws.Range("C103").Value = 42
$ws.Range("C103").NumberFormat = "0"

# --- DrawdownDuration (column C) cell removals: markers no longer present ---
$ws.Range("C52").Clear()
$ws.Range("C102").Clear()

# Touching ColumnWidth (read-only) forces the row span metadata back to
# "1:3" for the rows whose C cell was just cleared, matching the rest of the sheet.
$touch = $ws.Columns(3).ColumnWidth
